$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 16 new data rows after the existing last data row (row 23), copying the
# formatting of row 22 (the "interior" row style) so the final row (old row 23,
# now pushed down to row 39) keeps its special "last row" bottom-border style.
for ($i = 0; $i -lt 16; $i++) {
    $ws.Range("B23:J23").Insert()
    $ws.Range("B22:J22").Copy($ws.Range("B23:J23"))
}

# Fill in the 16 new employee rows (period 2509) in columns C (doc #), D (name),
# E (period), F (valor mora) and G (salario basico). B (Tipo Doc = "CC") is
# already correct from the copied row.
$ws.Range("C24").Value = "9153455"
$ws.Range("D24").Value = "HERNAN MEDINA NARVAEZ"
$ws.Range("E24").Value = "2509"
$ws.Range("F24").Value = 100000
$ws.Range("G24").Value = 2500000

$ws.Range("C25").Value = "1047429842"
$ws.Range("D25").Value = "ARTIS MANUEL ANAYA ROMERO"
$ws.Range("E25").Value = "2509"
$ws.Range("F25").Value = 68000
$ws.Range("G25").Value = 1700000

$ws.Range("C26").Value = "73101124"
$ws.Range("D26").Value = "EDGARDO ANTONIO VISBAL NARVAEZ"
$ws.Range("E26").Value = "2509"
$ws.Range("F26").Value = 120000
$ws.Range("G26").Value = 3000000

$ws.Range("C27").Value = "1143363823"
$ws.Range("D27").Value = "KAREN MARGARITA JULIO RODRIGUEZ"
$ws.Range("E27").Value = "2509"
$ws.Range("F27").Value = 56940
$ws.Range("G27").Value = 1423500

$ws.Range("C28").Value = "1002390885"
$ws.Range("D28").Value = "MARIA ALEJANDRA ARELLANO POLO"
$ws.Range("E28").Value = "2509"
$ws.Range("F28").Value = 56940
$ws.Range("G28").Value = 1423500

$ws.Range("C29").Value = "1001804871"
$ws.Range("D29").Value = "CATALINA SINNING BERTEL"
$ws.Range("E29").Value = "2509"
$ws.Range("F29").Value = 56940
$ws.Range("G29").Value = 1423500

$ws.Range("C30").Value = "1143412706"
$ws.Range("D30").Value = "NESTOR IVAN BUCHELI SANCHEZ"
$ws.Range("E30").Value = "2509"
$ws.Range("F30").Value = 56940
$ws.Range("G30").Value = 1423500

$ws.Range("C31").Value = "1143380806"
$ws.Range("D31").Value = "CESAR ANDRES ANGULO BARRIOS"
$ws.Range("E31").Value = "2509"
$ws.Range("F31").Value = 56940
$ws.Range("G31").Value = 1423500

$ws.Range("C32").Value = "1047512713"
$ws.Range("D32").Value = "YOSEP RAFAEL MOLINA CAMARGO"
$ws.Range("E32").Value = "2509"
$ws.Range("F32").Value = 68000
$ws.Range("G32").Value = 1700000

$ws.Range("C33").Value = "8772157"
$ws.Range("D33").Value = "ROBERTO ELIAS SOLANO MEJIA"
$ws.Range("E33").Value = "2509"
$ws.Range("F33").Value = 56940
$ws.Range("G33").Value = 1423500

$ws.Range("C34").Value = "1049928736"
$ws.Range("D34").Value = "NEIVY PAOLA PACHECO BERRIO"
$ws.Range("E34").Value = "2509"
$ws.Range("F34").Value = 72000
$ws.Range("G34").Value = 1800000

$ws.Range("C35").Value = "1140887488"
$ws.Range("D35").Value = "ENRIQUE ANTONIO PUELLO ROMERO"
$ws.Range("E35").Value = "2509"
$ws.Range("F35").Value = 56940
$ws.Range("G35").Value = 1423500

$ws.Range("C36").Value = "1001970138"
$ws.Range("D36").Value = "YULIS PATRICIA MORENO VILLALOBO"
$ws.Range("E36").Value = "2509"
$ws.Range("F36").Value = 56940
$ws.Range("G36").Value = 1423500

$ws.Range("C37").Value = "1043648823"
$ws.Range("D37").Value = "ANDREA CAROLINA VILLEGAS DIAZ"
$ws.Range("E37").Value = "2509"
$ws.Range("F37").Value = 56940
$ws.Range("G37").Value = 1423500

$ws.Range("C38").Value = "1004271922"
$ws.Range("D38").Value = "NEYDER YESID HERNANDEZ FONTALVO"
$ws.Range("E38").Value = "2509"
$ws.Range("F38").Value = 56940
$ws.Range("G38").Value = 1423500

$ws.Range("C39").Value = "1001897952"
$ws.Range("D39").Value = "AILIN PATRICIA MERCADO LIÃ?AN"
$ws.Range("E39").Value = "2509"
$ws.Range("F39").Value = 56940
$ws.Range("G39").Value = 1423500

# Update the summary figures: total "Valor Mora" in mora, worker count and period count.
$ws.Range("E11").Value = 1509860
$ws.Range("C13").Value = 16
$ws.Range("F13").Value = 9

# Column D widened to fit the new (longer) worker names.
$ws.Columns("D:D").ColumnWidth = 35.16
